$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 207
$ws.Range("F3").Value = 2485
$ws.Range("F5").Value = 1865
$ws.Range("F6").Value = 118
$ws.Range("F8").Value = 638
$ws.Range("F9").Value = 3642
$ws.Range("F15").Value = 1484
$ws.Range("F17").Value = 1813
$ws.Range("F20").Value = 15
$ws.Range("F22").Value = 1570
$ws.Range("F28").Value = 393
$ws.Range("F30").Value = 4411
$ws.Range("F31").Value = 75
$ws.Range("F32").Value = 75
$ws.Range("F37").Value = 1227
$ws.Range("F38").Value = 963

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F34").Value = 457

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F10").Value = 3041
$ws.Range("F11").Value = 550
$ws.Range("F12").Value = 841

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 207
$ws.Range("F5").Value = 2485
$ws.Range("F8").Value = 550
$ws.Range("F9").Value = 841
$ws.Range("F10").Value = 1866
$ws.Range("F14").Value = 1484
$ws.Range("F21").Value = 1813
$ws.Range("F23").Value = 15
$ws.Range("F26").Value = 1570
$ws.Range("F35").Value = 393
$ws.Range("F40").Value = 4411
$ws.Range("F41").Value = 75
$ws.Range("F42").Value = 457
$ws.Range("F51").Value = 1227
$ws.Range("F52").Value = 963
